$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header text (volume number, date range) ---
$ws.Range("A8").Value = "Volume 29   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/12/2022  Through  12/18/2022"

# --- Cells whose type/style changes (string <-> number) ---
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("F14").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 4
$ws.Range("K14").Copy($ws.Range("E16"))
$ws.Range("E16").Value = 25
$ws.Range("F14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("F14").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 1
$ws.Range("K14").Copy($ws.Range("E23"))
$ws.Range("E23").Value = 100
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("F14").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("K14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100
$ws.Range("C14").Copy($ws.Range("C29"))
$ws.Range("F14").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("K14").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100

# --- Simple value updates (same type/style, new number) ---
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("C16").Value = 5
$ws.Range("F16").Value = 26
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 136.363636363636
$ws.Range("I16").Value = 218
$ws.Range("J16").Value = 165
$ws.Range("K16").Value = 32.121212121212
$ws.Range("L16").Value = 80.165289256198
$ws.Range("M16").Value = -6.437768240343
$ws.Range("N16").Value = -50.454545454545
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -27.272727272727
$ws.Range("I17").Value = 254
$ws.Range("J17").Value = 226
$ws.Range("K17").Value = 12.389380530973
$ws.Range("L17").Value = 11.894273127753
$ws.Range("M17").Value = 44.318181818181
$ws.Range("N17").Value = 12.888888888888
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -26.666666666666
$ws.Range("I18").Value = 104
$ws.Range("J18").Value = 114
$ws.Range("K18").Value = -8.771929824561
$ws.Range("L18").Value = 6.122448979591
$ws.Range("M18").Value = -65.448504983388
$ws.Range("N18").Value = -86.717752234993
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 55.555555555555
$ws.Range("F19").Value = 68
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = 23.636363636363
$ws.Range("I19").Value = 583
$ws.Range("J19").Value = 485
$ws.Range("K19").Value = 20.206185567010
$ws.Range("L19").Value = 40.144230769230
$ws.Range("M19").Value = 23.516949152542
$ws.Range("N19").Value = 42.542787286063
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 80
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = 7.692307692307
$ws.Range("I20").Value = 297
$ws.Range("J20").Value = 304
$ws.Range("K20").Value = -2.302631578947
$ws.Range("L20").Value = 103.424657534247
$ws.Range("M20").Value = 45.588235294117
$ws.Range("N20").Value = -85.977337110481
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 25.925925925925
$ws.Range("F21").Value = 151
$ws.Range("G21").Value = 130
$ws.Range("H21").Value = 16.153846153846
$ws.Range("I21").Value = 1474
$ws.Range("J21").Value = 1313
$ws.Range("K21").Value = 12.261995430312
$ws.Range("L21").Value = 43.664717348927
$ws.Range("M21").Value = 4.836415362731
$ws.Range("N21").Value = -63.251059586138
$ws.Range("F22").Value = 2
$ws.Range("I22").Value = 16
$ws.Range("K22").Value = 23.076923076923
$ws.Range("L22").Value = 6.666666666666
$ws.Range("M22").Value = 23.076923076923
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = 300
$ws.Range("I23").Value = 34
$ws.Range("J23").Value = 44
$ws.Range("K23").Value = -22.727272727272
$ws.Range("L23").Value = -19.047619047619
$ws.Range("M23").Value = -33.333333333333
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 111.111111111111
$ws.Range("F24").Value = 129
$ws.Range("G24").Value = 92
$ws.Range("H24").Value = 40.217391304347
$ws.Range("I24").Value = 1160
$ws.Range("J24").Value = 975
$ws.Range("K24").Value = 18.974358974359
$ws.Range("L24").Value = 33.027522935779
$ws.Range("M24").Value = -16.486681065514
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 88.888888888888
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 14.285714285714
$ws.Range("I25").Value = 449
$ws.Range("J25").Value = 422
$ws.Range("K25").Value = 6.398104265402
$ws.Range("L25").Value = 19.098143236074
$ws.Range("M25").Value = 6.904761904761
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("D27").Value = 3
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -83.333333333333
$ws.Range("J27").Value = 44
$ws.Range("K27").Value = 6.818181818181
$ws.Range("L27").Value = 104.347826086957
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = -12.5
$ws.Range("J29").Value = 13
$ws.Range("K29").Value = -23.076923076923
